# Auto-generated: apply cryptos list price/volume updates (Tue Aug 27 2024 GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Force the cell to keep a literal text value (matches the source t="inlineStr" cells,
    # e.g. "550.17" or "8.20" must not be auto-coerced into the number 550.17 / 8.2).
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "62.016.38"
Set-TextValue "E2" "  -2.42%  "
Set-TextValue "D3" "2.577.79"
Set-TextValue "E3" "  -4.29%  "
Set-TextValue "D5" "550.17"
Set-TextValue "E5" "  -1.37%  "
Set-TextValue "D6" "155.37"
Set-TextValue "E6" "  -2.24%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "E8" "  +2.04%  "
Set-TextValue "E9" "  -1.20%  "
Set-TextValue "E10" "  -1.51%  "
Set-TextValue "D11" "5.49"
Set-TextValue "E11" "  +2.34%  "
Set-TextValue "E12" "  -0.91%  "
Set-TextValue "D13" "3.034.33"
Set-TextValue "E13" "  -4.33%  "
Set-TextValue "D14" "25.77"
Set-TextValue "E14" "  -2.84%  "
Set-TextValue "D15" "61.895.70"
Set-TextValue "E15" "  -2.38%  "
Set-TextValue "E16" "  +0.18%  "
Set-TextValue "D17" "2.581.15"
Set-TextValue "E17" "  -4.25%  "
Set-TextValue "D18" "11.64"
Set-TextValue "E18" "  -4.00%  "
Set-TextValue "E19" "  -0.32%  "
Set-TextValue "D20" "338.43"
Set-TextValue "E20" "  -2.17%  "
Set-TextValue "E21" "  -4.70%  "
Set-TextValue "E22" "  +0.34%  "
Set-TextValue "D23" "0.495"
Set-TextValue "E23" "  -2.69%  "
Set-TextValue "D24" "63.51"
Set-TextValue "E24" "  -0.76%  "
Set-TextValue "E25" "  -0.67%  "
Set-TextValue "E26" "  -0.05%  "
Set-TextValue "D27" "8.20"
Set-TextValue "E27" "  +0.20%  "
Set-TextValue "D28" "7.38"
Set-TextValue "E28" "  +4.61%  "
Set-TextValue "E29" "  -2.21%  "
Set-TextValue "E30" "  +0.61%  "
Set-TextValue "E31" "  -2.21%  "
Set-TextValue "D32" "162.54"
Set-TextValue "E32" "  -2.00%  "
Set-TextValue "D33" "4.88"
Set-TextValue "E33" "  +1.38%  "
Set-TextValue "E34" "  +0.03%  "
Set-TextValue "E35" "  +0.47%  "
Set-TextValue "E36" "  -1.82%  "
Set-TextValue "D37" "1.80"
Set-TextValue "E37" "  +0.44%  "
Set-TextValue "D38" "331.51"
Set-TextValue "D39" "6.05"
Set-TextValue "E39" "  -0.93%  "
Set-TextValue "D40" "0.918"
Set-TextValue "E40" "  -3.36%  "
Set-TextValue "E41" "  +0.08%  "
Set-TextValue "D42" "37.57"
Set-TextValue "E42" "  -1.75%  "
Set-TextValue "D43" "20.95"
Set-TextValue "E43" "  +0.10%  "
Set-TextValue "E44" "  +0.00%  "
Set-TextValue "E45" "  -2.61%  "
Set-TextValue "D46" "2.120.25"
Set-TextValue "E46" "  +0.48%  "
Set-TextValue "D47" "0.0550"
Set-TextValue "E47" "  -2.90%  "
Set-TextValue "E48" "  -1.16%  "
Set-TextValue "D49" "19.61"
Set-TextValue "E49" "  -4.09%  "
Set-TextValue "E50" "  -0.62%  "
Set-TextValue "E51" "  -0.91%  "
